# Add a new "VAT" worksheet (positioned after "Contacts") listing each
# seller/supplier and whether they are VAT registered, with a hyperlink to
# each seller's website/profile where one is known.

$wb = $excel.ActiveWorkbook
$contacts = $wb.Worksheets.Item("Contacts")

$vat = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $contacts)
$vat.Name = "VAT"

# Seller / VAT-registered rows (row 11 - Brighthampton Veg - has no known website)
$sellers = @(
    @("Blacklands", "N", "https://www.facebook.com/BlacklandsOrganics/?locale=en_GB"),
    @("Worthy Earth at Bletchingdon", "N", "https://worthy-earth.co.uk/"),
    @("Bruern Farms", "Y", "https://bruernfarms.co.uk/"),
    @("Nettlebed Creamery", "Y", "http://www.nettlebedcreamery.com/"),
    @("Oxford City Farm", "N", "http://www.oxfordcityfarm.org.uk/"),
    @("Rose Dale's Organic Farm", "N", "http://ww.rosedalesorganicfarm.com/"),
    @("Sandy Lane Organics", "N", "http://sandylanefarm.co.uk/"),
    @("Willowbrook Farm", "N", "https://www.willowbrookfarm.co.uk/"),
    @("Forge Farm", "N", "https://www.forgefarm.com/"),
    @("Brighthampton Veg", "N", ""),
    @("Pudlicote Farm", "N", "www.pudlicotefarm.co.uk")
)

# Add the hyperlinks first (this stamps the cell with the address as its
# text); the real seller name written afterwards then overwrites the cell
# value while the hyperlink keeps its own cached display text, and
# resetting the style back to Normal drops the blue/underline formatting
# that Hyperlinks.Add applies by default.
for ($i = 0; $i -lt $sellers.Count; $i++) {
    $row = $i + 2
    $url = $sellers[$i][2]
    if ($url -ne "") {
        $vat.Hyperlinks.Add($vat.Range("A$row"), $url, "", "", $url)
    }
}

# Header row
$vat.Range("A1").Value = "Seller"
$vat.Range("B1").Value = "VAT Y/N?"

for ($i = 0; $i -lt $sellers.Count; $i++) {
    $row = $i + 2
    $vat.Range("A$row").Value = $sellers[$i][0]
    $vat.Range("B$row").Value = $sellers[$i][1]
    if ($sellers[$i][2] -ne "") {
        $vat.Range("A$row").Style = "Normal"
    }
}

# Column widths (A wide for seller names, B narrow for the Y/N flag)
$vat.Columns.Item(1).ColumnWidth = 27.14
$vat.Columns.Item(2).ColumnWidth = 8.92

# Match the authored selection/active-cell state
$vat.Range("A8").Select() | Out-Null
